# Update TPM-derived ligand-receptor edge metrics (Tg -> Asgr1)
# for Young D7 lrc2p NATMI output, per new TPM script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07266366666666667
$ws.Range("H2").Value = 0.217991
$ws.Range("I2").Value = 0.08800019376989421
$ws.Range("J2").Value = 0.08800019376989422
$ws.Range("K2").Value = 2.0
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5478816666666667
$ws.Range("N2").Value = 1.643645
$ws.Range("O2").Value = 0.1406795382009894
$ws.Range("P2").Value = 0.1406795382009894
$ws.Range("Q2").Value = 0.03981109079944445
$ws.Range("R2").Value = 0.358299817195
$ws.Range("S2").Value = 0.0123798266211463
$ws.Range("T2").Value = 0.0123798266211463

# Row 3
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07266366666666667
$ws.Range("H3").Value = 0.217991
$ws.Range("I3").Value = 0.08800019376989421
$ws.Range("J3").Value = 0.08800019376989422
$ws.Range("K3").Value = 3.0
$ws.Range("O3").Value = 0.7256059791788686
$ws.Range("P3").Value = 0.7256059791788687
$ws.Range("Q3").Value = 0.2053402071908889
$ws.Range("R3").Value = 1.848061864718
$ws.Range("S3").Value = 0.06385346676833425
$ws.Range("T3").Value = 0.06385346676833428

# Row 4
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07266366666666667
$ws.Range("H4").Value = 0.217991
$ws.Range("I4").Value = 0.08800019376989421
$ws.Range("J4").Value = 0.08800019376989422
$ws.Range("K4").Value = 2.0
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4691703333333333
$ws.Range("N4").Value = 1.407511
$ws.Range("O4").Value = 0.1204688345067291
$ws.Range("P4").Value = 0.1204688345067291
$ws.Range("Q4").Value = 0.03409163671122222
$ws.Range("R4").Value = 0.306824730401
$ws.Range("S4").Value = 0.01060128077982548
$ws.Range("T4").Value = 0.01060128077982548

# Row 5
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07266366666666667
$ws.Range("H5").Value = 0.217991
$ws.Range("I5").Value = 0.08800019376989421
$ws.Range("J5").Value = 0.08800019376989422
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05158566666666667
$ws.Range("N5").Value = 0.154757
$ws.Range("O5").Value = 0.01324564811341288
$ws.Range("P5").Value = 0.01324564811341288
$ws.Range("Q5").Value = 0.003748403687444445
$ws.Range("R5").Value = 0.033735633187
$ws.Range("S5").Value = 0.001165619600588167
$ws.Range("T5").Value = 0.001165619600588168

# Row 6
$ws.Range("I6").Value = 0.1096079590984048
$ws.Range("J6").Value = 0.1096079590984048
$ws.Range("K6").Value = 2.0
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5478816666666667
$ws.Range("N6").Value = 1.643645
$ws.Range("O6").Value = 0.1406795382009894
$ws.Range("P6").Value = 0.1406795382009894
$ws.Range("Q6").Value = 0.04958639549611112
$ws.Range("R6").Value = 0.446277559465
$ws.Range("S6").Value = 0.01541959706911653
$ws.Range("T6").Value = 0.01541959706911653

# Row 7
$ws.Range("I7").Value = 0.1096079590984048
$ws.Range("J7").Value = 0.1096079590984048
$ws.Range("K7").Value = 3.0
$ws.Range("O7").Value = 0.7256059791788686
$ws.Range("P7").Value = 0.7256059791788687
$ws.Range("Q7").Value = 0.2557599030962223
$ws.Range("S7").Value = 0.07953219048739542
$ws.Range("T7").Value = 0.07953219048739543

# Row 8
$ws.Range("I8").Value = 0.1096079590984048
$ws.Range("J8").Value = 0.1096079590984048
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.4691703333333333
$ws.Range("N8").Value = 1.407511
$ws.Range("O8").Value = 0.1204688345067291
$ws.Range("P8").Value = 0.1204688345067291
$ws.Range("Q8").Value = 0.04246257379855555
$ws.Range("R8").Value = 0.382163164187
$ws.Range("S8").Value = 0.01320434308524607
$ws.Range("T8").Value = 0.01320434308524607

# Row 9
$ws.Range("I9").Value = 0.1096079590984048
$ws.Range("J9").Value = 0.1096079590984048
$ws.Range("K9").Value = 1.0
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05158566666666667
$ws.Range("N9").Value = 0.154757
$ws.Range("O9").Value = 0.01324564811341288
$ws.Range("P9").Value = 0.01324564811341288
$ws.Range("Q9").Value = 0.004668795152111111
$ws.Range("R9").Value = 0.04201915636900001
$ws.Range("S9").Value = 0.001451828456646822
$ws.Range("T9").Value = 0.001451828456646823

# Row 10
$ws.Range("G10").Value = 0.5710436666666667
$ws.Range("H10").Value = 1.713131
$ws.Range("I10").Value = 0.6915691930089437
$ws.Range("J10").Value = 0.6915691930089438
$ws.Range("K10").Value = 2.0
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.5478816666666667
$ws.Range("N10").Value = 1.643645
$ws.Range("O10").Value = 0.1406795382009894
$ws.Range("P10").Value = 0.1406795382009894
$ws.Range("Q10").Value = 0.3128643558327778
$ws.Range("R10").Value = 2.815779202495
$ws.Range("S10").Value = 0.09728963470652911
$ws.Range("T10").Value = 0.09728963470652913

# Row 11
$ws.Range("G11").Value = 0.5710436666666667
$ws.Range("H11").Value = 1.713131
$ws.Range("I11").Value = 0.6915691930089437
$ws.Range("J11").Value = 0.6915691930089438
$ws.Range("K11").Value = 3.0
$ws.Range("O11").Value = 0.7256059791788686
$ws.Range("P11").Value = 0.7256059791788687
$ws.Range("Q11").Value = 1.613711916937556
$ws.Range("R11").Value = 14.523407252438
$ws.Range("S11").Value = 0.5018067414631946
$ws.Range("T11").Value = 0.5018067414631947

# Row 12
$ws.Range("G12").Value = 0.5710436666666667
$ws.Range("H12").Value = 1.713131
$ws.Range("I12").Value = 0.6915691930089437
$ws.Range("J12").Value = 0.6915691930089438
$ws.Range("K12").Value = 2.0
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.4691703333333333
$ws.Range("N12").Value = 1.407511
$ws.Range("O12").Value = 0.1204688345067291
$ws.Range("P12").Value = 0.1204688345067291
$ws.Range("Q12").Value = 0.2679167474378889
$ws.Range("R12").Value = 2.411250726941
$ws.Range("S12").Value = 0.08331253466254665
$ws.Range("T12").Value = 0.08331253466254666

# Row 13
$ws.Range("G13").Value = 0.5710436666666667
$ws.Range("H13").Value = 1.713131
$ws.Range("I13").Value = 0.6915691930089437
$ws.Range("J13").Value = 0.6915691930089438
$ws.Range("K13").Value = 1.0
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05158566666666667
$ws.Range("N13").Value = 0.154757
$ws.Range("O13").Value = 0.01324564811341288
$ws.Range("P13").Value = 0.01324564811341288
$ws.Range("Q13").Value = 0.02945766824077778
$ws.Range("R13").Value = 0.265119014167
$ws.Range("S13").Value = 0.009160282176673385
$ws.Range("T13").Value = 0.009160282176673386

# Row 14
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.09150866666666667
$ws.Range("H14").Value = 0.274526
$ws.Range("I14").Value = 0.1108226541227573
$ws.Range("J14").Value = 0.1108226541227573
$ws.Range("K14").Value = 2.0
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.5478816666666667
$ws.Range("N14").Value = 1.643645
$ws.Range("O14").Value = 0.1406795382009894
$ws.Range("P14").Value = 0.1406795382009894
$ws.Range("Q14").Value = 0.05013592080777778
$ws.Range("R14").Value = 0.45122328727
$ws.Range("S14").Value = 0.01559047980419747
$ws.Range("T14").Value = 0.01559047980419747

# Row 15
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.09150866666666667
$ws.Range("H15").Value = 0.274526
$ws.Range("I15").Value = 0.1108226541227573
$ws.Range("J15").Value = 0.1108226541227573
$ws.Range("K15").Value = 3.0
$ws.Range("O15").Value = 0.7256059791788686
$ws.Range("P15").Value = 0.7256059791788687
$ws.Range("Q15").Value = 0.2585942801275556
$ws.Range("R15").Value = 2.327348521148
$ws.Range("S15").Value = 0.08041358045994437
$ws.Range("T15").Value = 0.08041358045994439

# Row 16
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.09150866666666667
$ws.Range("H16").Value = 0.274526
$ws.Range("I16").Value = 0.1108226541227573
$ws.Range("J16").Value = 0.1108226541227573
$ws.Range("K16").Value = 2.0
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.4691703333333333
$ws.Range("N16").Value = 1.407511
$ws.Range("O16").Value = 0.1204688345067291
$ws.Range("P16").Value = 0.1204688345067291
$ws.Range("Q16").Value = 0.04293315164288888
$ws.Range("R16").Value = 0.386398364786
$ws.Range("S16").Value = 0.01335067597911093
$ws.Range("T16").Value = 0.01335067597911093

# Row 17
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.09150866666666667
$ws.Range("H17").Value = 0.274526
$ws.Range("I17").Value = 0.1108226541227573
$ws.Range("J17").Value = 0.1108226541227573
$ws.Range("K17").Value = 1.0
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.05158566666666667
$ws.Range("N17").Value = 0.154757
$ws.Range("O17").Value = 0.01324564811341288
$ws.Range("P17").Value = 0.01324564811341288
$ws.Range("Q17").Value = 0.004720535575777778
$ws.Range("R17").Value = 0.042484820182
$ws.Range("S17").Value = 0.001467917879504508
$ws.Range("T17").Value = 0.001467917879504508
